$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B width (Excel column width units) to match target (15.4 -> 14.3)
$ws.Columns.Item(2).ColumnWidth = 14.3

# Row 6 updated values (B6:F6)
$ws.Range("B6").Value = -149661000.0
$ws.Range("C6").Value = -160289000.0
$ws.Range("D6").Value = -118497000.0
$ws.Range("E6").Value = -52512000.0
$ws.Range("F6").Value = -22544000.0

# Row 8 updated values (B8:F8)
$ws.Range("B8").Value = 581000000.0
$ws.Range("C8").Value = 576000000.0
$ws.Range("D8").Value = 445515000.0
$ws.Range("E8").Value = 307744000.0
$ws.Range("F8").Value = 172100000.0
